# Adds "start DD" / "start SD" label cells and corrects the walk/no-walk
# labels in columns C and I for the ten_lists sheet (OA S019), and renames
# the S000 shared string to S019.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ten_lists")

# Rename the sheet-name placeholder string used elsewhere (S000 -> S019)
$found = $ws.Cells.Find("S000")
if ($found -ne $null) {
    $found.Value = "S019"
}

# Row 3
$ws.Range("E3").Value = "start DD"
$ws.Range("K3").Value = "start DD"

# Row 10
$ws.Range("C10").Value = "no walk/same"
$ws.Range("E10").Value = "start SD"
$ws.Range("I10").Value = "no walk/diff"
$ws.Range("K10").Value = "start DD"

# Row 17
$ws.Range("C17").Value = "walk/diff"
$ws.Range("E17").Value = "start SD"
$ws.Range("I17").Value = "walk/diff"
$ws.Range("K17").Value = "start SD"

# Row 24
$ws.Range("C24").Value = "no walk/same"
$ws.Range("E24").Value = "start DD"
$ws.Range("I24").Value = "no walk/same"
$ws.Range("K24").Value = "start SD"

# Row 31
$ws.Range("C31").Value = "no walk/diff"
$ws.Range("E31").Value = "start SD"
$ws.Range("I31").Value = "walk/same"
$ws.Range("K31").Value = "start DD"

# Best-fit the two newly-relevant columns (A holds the row numbers 1-5,
# F holds the sheet code) to match the narrower widths Excel would compute.
$ws.Columns("A:A").ColumnWidth = 0.92
$ws.Columns("F:F").ColumnWidth = 3.75

$ws.Range("K25").Select()
